# Penalty/Reward system update: refresh the forecast comparison table and
# the summary statistics that are derived from it.
#
# Note: several "Value"-looking strings (dates, and numbers-as-text in the
# Summary sheet) must stay text cells, not get auto-converted by Excel into
# dates/numbers. We use the classic leading-apostrophe (quote-prefix) trick
# for those so they are stored as literal text, exactly like the source data.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Each week's Week_Start_Date (col B) rolls forward by one week, and the
# MyForecast value (col D) is refreshed with the latest forecast numbers.

$forecastRows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 227 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 304 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 373 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 420 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 296 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 297 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 280 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 290 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 283 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 283 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 278 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 277 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 274 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 362 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 310 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 301 }
)

foreach ($entry in $forecastRows) {
    $wsForecast.Cells.Item($entry.Row, 2).Value = "'" + $entry.Date
    $wsForecast.Cells.Item($entry.Row, 4).Value = $entry.Forecast
}

# --- Summary sheet ---------------------------------------------------------
# Recomputed stats to match the refreshed forecast/historical window.

$wsSummary.Range("B2").Value  = "'2023-02-26 to 2025-01-05"
$wsSummary.Range("B4").Value  = "'410"
$wsSummary.Range("B5").Value  = "'199"
$wsSummary.Range("B6").Value  = "'192"
$wsSummary.Range("B7").Value  = "'113"
$wsSummary.Range("B8").Value  = "'12132 units"
$wsSummary.Range("B9").Value  = "'4855"
$wsSummary.Range("B10").Value = "'2486"
$wsSummary.Range("B11").Value = "'1324"
$wsSummary.Range("B12").Value = "'420"
$wsSummary.Range("B14").Value = "'227"
$wsSummary.Range("B15").Value = "'2025-01-12"
